$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MACRO_SCORE (column N) for rows 2-6 with the refreshed value
$newValue = 85.92117485762657

$ws.Range("N2:N6").Value = $newValue
